# Updates cryptos list values (Price / Volume(1h) columns) per the
# 2024-05-03 GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. '577.11') must be forced to stay plain text, matching the original
# inlineStr cells, and then have their style reset to 'Normal' so no stray
# number-format style lingers on the cell.
$numericLookingCells = @{
    'D5' = '577.11'
    'D6' = '143.14'
    'D11' = '5.49'
    'D12' = '0.465'
    'D14' = '34.94'
    'D20' = '449.64'
    'D21' = '13.96'
    'D22' = '0.733'
    'D23' = '7.31'
    'D24' = '13.58'
    'D25' = '81.86'
    'D27' = '2.23'
    'D30' = '8.06'
    'D31' = '6.66'
    'D32' = '26.61'
    'D36' = '6.06'
    'D37' = '2.19'
    'D38' = '50.12'
    'D39' = '2.97'
    'D40' = '8.81'
    'D41' = '419.10'
    'D42' = '0.0368'
    'D44' = '0.109'
    'D46' = '36.58'
    'D47' = '2.10'
    'D49' = '123.27'
    'D51' = '24.10'
}
foreach ($addr in $numericLookingCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingCells[$addr]
    $cell.Style = "Normal"
}

# Remaining cells are safe to assign directly: they already contain
# non-numeric characters (extra '.' grouping separators, '%' signs, padding
# spaces, or the subscript digit in the PEPE price) so Excel keeps them as text.
$textCells = @{
    'D2' = '61.906.88'
    'E2' = '  +4.54%  '
    'D3' = '3.064.64'
    'E3' = '  +2.30%  '
    'E4' = '  -0.09%  '
    'E5' = '  +2.64%  '
    'E6' = '  +3.33%  '
    'E7' = '  +0.07%  '
    'D8' = '3.051.67'
    'E8' = '  +2.15%  '
    'E9' = '  +0.85%  '
    'E10' = '  +4.77%  '
    'E11' = '  +11.57%  '
    'E12' = '  +1.85%  '
    'E13' = '  +4.32%  '
    'E14' = '  +3.45%  '
    'E15' = '  -0.01%  '
    'D16' = '3.569.08'
    'E16' = '  +2.26%  '
    'E17' = '  +3.18%  '
    'D18' = '3.058.86'
    'E18' = '  +2.47%  '
    'D19' = '61.819.62'
    'E19' = '  +4.37%  '
    'E20' = '  +5.70%  '
    'E21' = '  +3.15%  '
    'E22' = '  +3.07%  '
    'E23' = '  +2.28%  '
    'E24' = '  +1.09%  '
    'E25' = '  +1.67%  '
    'E26' = '  +0.09%  '
    'E27' = '  +4.65%  '
    'E28' = '  -0.29%  '
    'E29' = '  +3.98%  '
    'E30' = '  +3.17%  '
    'E31' = '  +8.31%  '
    'E32' = '  +3.78%  '
    'E33' = '  +7.64%  '
    'D34' = '0.0₃0808'
    'E34' = '  +4.79%  '
    'E35' = '  +2.11%  '
    'E36' = '  +5.44%  '
    'E37' = '  +4.86%  '
    'E38' = '  +2.19%  '
    'E39' = '  +5.77%  '
    'E40' = '  +1.78%  '
    'E41' = '  +4.07%  '
    'E42' = '  +5.14%  '
    'D43' = '2.774.40'
    'E43' = '  +1.11%  '
    'E44' = '  +1.04%  '
    'E45' = '  +7.78%  '
    'E46' = '  +12.18%  '
    'E47' = '  +3.14%  '
    'E48' = '  -0.03%  '
    'E49' = '  -1.69%  '
    'E50' = '  +1.57%  '
    'E51' = '  +2.65%  '
}
foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}

